# Apply "nuevos experimentos no convexos" value updates.
#
# All target cells hold text (shared-string) values, even the ones that
# look numeric (e.g. "0.77", "4.5"). Writing a numeric-looking string via
# .Value would make Excel auto-convert it to a real number, so we prefix
# the value with a leading apostrophe to force text entry, then reset the
# cell style back to "Normal" so no stray quote-prefix style is left
# behind on the cell.

function Set-TextValue {
    param(
        $Worksheet,
        [string]$Address,
        [string]$Text
    )
    $Worksheet.Range($Address).Value = "'" + $Text
    $Worksheet.Range($Address).Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# --- Restricciones_del_follower ---------------------------------------
$wsFollower = $wb.Worksheets.Item("Restricciones_del_follower")

Set-TextValue $wsFollower "A2" "1.1000000000000056 - 2x_1 + y_1 - y_2"
Set-TextValue $wsFollower "B2" "1.3999999999999944"
Set-TextValue $wsFollower "D2" "0.92"
Set-TextValue $wsFollower "F2" "0.8999999999999999"

Set-TextValue $wsFollower "A3" "2.5499999999999927 + x_1 - 3x_2 + y_2"
Set-TextValue $wsFollower "B3" "-4.549999999999993"
Set-TextValue $wsFollower "D3" "0.36"
Set-TextValue $wsFollower "E3" "6.5"

Set-TextValue $wsFollower "A4" "102.2 - y_1"
Set-TextValue $wsFollower "B4" "-102.2"
Set-TextValue $wsFollower "D4" "0.49"
Set-TextValue $wsFollower "F4" "2.2"

Set-TextValue $wsFollower "A5" "-0.2 - y_2"
Set-TextValue $wsFollower "B5" "-0.2"
Set-TextValue $wsFollower "D5" "0.75"
Set-TextValue $wsFollower "F5" "1.2"

# --- Punto_modificado ---------------------------------------------------
$wsPunto = $wb.Worksheets.Item("Punto_modificado")

Set-TextValue $wsPunto "A2" "51.550000000000004"
Set-TextValue $wsPunto "B2" "18.099999999999998"
Set-TextValue $wsPunto "C2" "102.2"
Set-TextValue $wsPunto "D2" "0.2"

# --- Vector_bf -----------------------------------------------------------
# Use numeric index (5) here too, to stay consistent with the "Vector_BF"
# sheet below and avoid any ambiguity from the case-insensitive name match
# against the similarly-named "Vector_BF" sheet.
$wsBf = $wb.Worksheets.Item(5)

Set-TextValue $wsBf "A2" "3.57"
Set-TextValue $wsBf "A3" "0.31000000000000005"

# --- Vector_BF -----------------------------------------------------------
# NOTE: worksheet name lookup by string is case-insensitive, and this
# workbook has both "Vector_bf" and "Vector_BF" sheets whose names differ
# only by case. Use the sheet's numeric index (6) to unambiguously target
# "Vector_BF" instead of relying on name resolution.
$wsBF = $wb.Worksheets.Item(6)

Set-TextValue $wsBF "A2" "-4.5"
Set-TextValue $wsBF "A3" "18.5"
Set-TextValue $wsBF "A5" "-6.5"
